# Updated assumptions to match class diagram
#
# Appends a series of new "List Paragraph" (numId 1) bullet items after the
# existing "There is a separate entity ..." bullet at the end of the
# Assumptions section.

$d = $word.ActiveDocument

function Add-ListParagraph($RunTexts) {
    # Create a brand-new (empty) paragraph right after the current last
    # paragraph in the document, inheriting its paragraph style/numbering.
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $nr = $newPara.Range

    # Build the run markup for this paragraph, preserving each run
    # individually (matching the source structure) and adding
    # xml:space="preserve" whenever leading/trailing whitespace matters.
    $runsXml = ""
    foreach ($t in $RunTexts) {
        $escaped = $t -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
        $preserve = ""
        if ($t -ne $t.Trim() -or $t -eq "") {
            $preserve = ' xml:space="preserve"'
        }
        $runsXml += "<w:r><w:t$preserve>$escaped</w:t></w:r>"
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $nr.InsertXML($xml)

    # InsertXML leaves an extra empty trailing paragraph behind (it inserts
    # the new paragraph mark before the one that was already there). Merge
    # that stray empty paragraph away by deleting across the boundary
    # between the new content paragraph and it.
    $contentPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $delRange = $d.Range($contentPara.Range.End - 1, $trailingPara.Range.End)
    $delRange.Delete()
}

Add-ListParagraph @("There is no call centre operator")
Add-ListParagraph @("Call distribution ", "is handled by an automated system")
Add-ListParagraph @("RMs and HR managers are both types of employees")
Add-ListParagraph @("A HR manager ", "oversees at least ", "one or more RMs")
Add-ListParagraph @("A RM ", "must have", " one HR manager")
Add-ListParagraph @("A RM ", "can serve many customers")
Add-ListParagraph @("A customer is only served by one RM")
Add-ListParagraph @("A customer can make", " many", " orders")
Add-ListParagraph @("An order can only be made by one customer")
Add-ListParagraph @("An order can be for one or many travel packages")
Add-ListParagraph @("An order can only have one payment")
Add-ListParagraph @("A payment ", "can be made for multiple orders")
